# packaging co2mpas driver model
#
# 1) "inputs" sheet: rename the input row from "v_start" to "starting_speed"
#    (label in column A), leaving the "Current speed" description in column C
#    untouched; leave the active cell on C14.
# 2) "time_series" sheet: replace the literal "velocities" column (B) with a
#    running-count formula column (A) -- A3 = A2+1, filled down through A22 --
#    and clear out column B entirely; leave the active cell on O8.
#
# Sheet "time_series" must remain the active tab when we're done (it was the
# active tab originally), so touch "inputs" first and finish on "time_series".

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("inputs")
$wsInputs.Range("A4").Value = "starting_speed"
$wsInputs.Range("C14").Select()

$wsTime = $wb.Worksheets.Item("time_series")
$wsTime.Activate()
$wsTime.Range("B2:B22").ClearContents()
$wsTime.Range("A3").Formula = "=A2+1"
$wsTime.Range("A4:A22").Formula = "=A3+1"
$wsTime.Range("O8").Select()
